# Apply the coin-price / volume snapshot refresh described by the commit.
#
# All data cells on this sheet are stored as *text* (prices use "."
# thousand separators like "61.964.36", volumes are "  +0.57%  " with
# padding) even though many of the price strings would otherwise parse as
# plain numbers. Writing such a string straight into `.Value` lets Excel
# auto-coerce it to a Double, which would change the stored cell type.
# Flipping the cell to the "@" (Text) number format first forces the
# literal text to stick, and re-applying the "Normal" style right after
# drops that temporary formatting again so the cell ends up styled exactly
# like it started (no stray NumberFormat left behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" '61.964.36'
Set-TextValue "E2" '  -0.04%  '
Set-TextValue "D3" '3.434.96'
Set-TextValue "E3" '  +0.57%  '
Set-TextValue "E4" '  +0.10%  '
Set-TextValue "D5" '412.84'
Set-TextValue "E5" '  +0.62%  '
Set-TextValue "D6" '129.72'
Set-TextValue "E6" '  +0.85%  '
Set-TextValue "E7" '  +1.76%  '
Set-TextValue "E8" '  +0.10%  '
Set-TextValue "D9" '0.733'
Set-TextValue "E9" '  -2.21%  '
Set-TextValue "D10" '0.142'
Set-TextValue "E10" '  +1.21%  '
Set-TextValue "D11" '43.55'
Set-TextValue "E11" '  +1.60%  '
Set-TextValue "D12" '0.0000219'
Set-TextValue "E12" '  +9.39%  '
Set-TextValue "D13" '9.27'
Set-TextValue "E13" '  +4.79%  '
Set-TextValue "D14" '3.981.13'
Set-TextValue "E14" '  +0.60%  '
Set-TextValue "E15" '  +0.47%  '
Set-TextValue "D16" '21.17'
Set-TextValue "E16" '  -0.52%  '
Set-TextValue "D17" '3.427.28'
Set-TextValue "E17" '  +0.61%  '
Set-TextValue "D18" '12.67'
Set-TextValue "E18" '  -0.16%  '
Set-TextValue "E19" '  +3.61%  '
Set-TextValue "D20" '61.999.47'
Set-TextValue "D21" '487.72'
Set-TextValue "E21" '  +19.99%  '
Set-TextValue "D22" '92.71'
Set-TextValue "E22" '  +2.23%  '
Set-TextValue "E23" '  +3.84%  '
Set-TextValue "D24" '13.60'
Set-TextValue "E24" '  +1.37%  '
Set-TextValue "D25" '3.39'
Set-TextValue "E25" '  +4.46%  '
Set-TextValue "D26" '34.77'
Set-TextValue "E26" '  +4.94%  '
Set-TextValue "D27" '9.14'
Set-TextValue "E27" '  +6.98%  '
Set-TextValue "D28" '4.81'
Set-TextValue "E28" '  +0.36%  '
Set-TextValue "D29" '7.69'
Set-TextValue "E29" '  +0.69%  '
Set-TextValue "E30" '  -0.95%  '
Set-TextValue "E31" '  +2.78%  '
Set-TextValue "B32" 'Kaspa'
Set-TextValue "C32" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D32" '0.169'
Set-TextValue "E32" '  -1.86%  '
Set-TextValue "B33" 'Hedera'
Set-TextValue "C33" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D33" '0.114'
Set-TextValue "E33" '  -2.50%  '
Set-TextValue "D34" '42.07'
Set-TextValue "E34" '  -4.12%  '
Set-TextValue "E35" '  +0.04%  '
Set-TextValue "D36" '58.69'
Set-TextValue "E36" '  +11.35%  '
Set-TextValue "D37" '0.0497'
Set-TextValue "E37" '  -0.31%  '
Set-TextValue "D38" '3.50'
Set-TextValue "E38" '  +3.06%  '
Set-TextValue "E39" '  +0.09%  '
Set-TextValue "D40" '150.93'
Set-TextValue "E40" '  +7.42%  '
Set-TextValue "D41" '2.18'
Set-TextValue "E41" '  +10.39%  '
Set-TextValue "E42" '  +4.09%  '
Set-TextValue "E43" '  +3.42%  '
Set-TextValue "E44" '  +2.42%  '
Set-TextValue "D45" '2.67'
Set-TextValue "E45" '  +12.63%  '
Set-TextValue "E46" '  +6.88%  '
Set-TextValue "D47" '2.41'
Set-TextValue "E47" '  +25.56%  '
Set-TextValue "D48" '16.64'
Set-TextValue "E48" '  -0.81%  '
Set-TextValue "D49" '23.00'
Set-TextValue "E49" '  +5.82%  '
Set-TextValue "D50" '118.61'
Set-TextValue "E50" '  +22.89%  '
Set-TextValue "D51" '0.146'
Set-TextValue "E51" '  +16.04%  '
